$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Sex" column (J) ---------------------------------------------
# Header
$ws.Range("J1").Value = "Sex"

# Rows 2 & 3 pick up the plain (unformatted) column default style, same
# as the existing untouched cells in those rows.
$ws.Range("J2").Value = "Male"
$ws.Range("J3").Value = "Male"

# Rows 4-7 copy the date-column (I) number format/style onto J before
# writing the text value, matching style index 4 used by I4:I7.
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("J4").Value = "Female"

$ws.Range("I5").Copy()
$ws.Range("J5").PasteSpecial(-4122)
$ws.Range("J5").Value = "Male"

$ws.Range("I6").Copy()
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("J6").Value = "Male"

$ws.Range("I7").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$ws.Range("J7").Value = "Female"

# --- Backfill the Date of Birth formulas that were missing -------------
$ws.Range("I4").Formula = "=DATE(2004,9,22)"
$ws.Range("I5").Formula = "=DATE(2005,1,19)"
$ws.Range("I6").Formula = "=DATE(2004,6,20)"
$ws.Range("I7").Formula = "=DATE(2003,7,30)"

# --- Column J width (best-fit sized for "Female") -----------------------
$ws.Columns.Item(10).ColumnWidth = 9.5

# --- Update the saved cursor/selection ----------------------------------
$ws.Range("I15").Select()
